$wb = $excel.ActiveWorkbook

# --- 1. Update the "Last Updated" timestamp on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(2, 1).Value = "05 Nov 2025, 11:08 AM"

# --- 2. Refresh the "Stock List" sheet ---
# The market-health refresh dropped the old top row (CAPTRU-RE1), so every
# remaining stock shifts up by one row, and a brand-new stock (TRAVELFOOD)
# is appended as the new last row.
$ws = $wb.Worksheets.Item("Stock List")

# Remove the old first data row; rows below shift up automatically.
$ws.Rows.Item(2).Delete()

# Append the new row of data at the bottom (row 76).
$lastRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($lastRow, 1).Value = [char]0x1F4CB
$ws.Cells.Item($lastRow, 2).Value = "TRAVELFOOD"
$ws.Cells.Item($lastRow, 3).Value = "TRAVELFOOD"
$ws.Cells.Item($lastRow, 4).Value = 1316.3
$ws.Cells.Item($lastRow, 5).Value = 0.1141
$ws.Cells.Item($lastRow, 6).Value = "N/A"
$ws.Cells.Item($lastRow, 7).Value = "N/A"
$ws.Cells.Item($lastRow, 8).Value = 17332.9705
